$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("template")
$ws2 = $wb.Worksheets.Item("example")

# 1) template!A5 = "Group" (first new shared string -> index 9)
$ws.Range("A1").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "Group"

# 2) example!A4 = "Group" (reuses index 9), example!B4 = "Site" (new -> index 10)
$ws2.Range("A1").Copy()
$ws2.Range("A4").PasteSpecial(-4122)
$ws2.Range("A4").Value = "Group"

$ws2.Range("B1").Copy()
$ws2.Range("B4").PasteSpecial(-4122)
$ws2.Range("B4").Value = "Site"

# 3) template!B5 = "Group category..." (new -> index 11)
$ws.Range("B1").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = "Group category used to break data into different plots (e.g. site, provider, patient demographics, etc.)  "

# Update selections to match the new active cells on each sheet
$ws2.Activate()
$ws2.Range("A4:B4").Select()

$ws.Activate()
$ws.Range("D14").Select()
